$wb = $excel.ActiveWorkbook

# Update the Environment value in Constants!B19 from DEV to PROD.
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("B19").Value = "PROD"

# Make Constants the active sheet (tab selected) as captured in the diff
# (activeTab moves from Assets to Constants).
$constants.Activate()
